$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: PV_DAO_02 - "Kiem tra ton tai (checkExist)" ---
$ws.Range("A2").Value = "PV_DAO_02"
$ws.Range("B2").Value = "Kiểm tra tồn tại (checkExist)"
$ws.Range("C2").Value = "Item đã tồn tại"
$ws.Range("D2").Value = "Insert ProductId=1, SizeId=14"
$ws.Range("E2").Value = "Return Object (ID, Stock)"

# --- Row 3: PV_DAO_03 - "Cap nhat ton kho" ---
$ws.Range("A3").Value = "PV_DAO_03"
$ws.Range("B3").Value = "Cập nhật tồn kho"
$ws.Range("C3").Value = "New Stock=9999"
$ws.Range("D3").Value = "Update Stock 9999"
$ws.Range("E3").Value = "Stock = 9999"

# --- Row 4: PV_DAO_04 - "Tim bien the theo SP" ---
$ws.Range("A4").Value = "PV_DAO_04"
$ws.Range("B4").Value = "Tìm biến thể theo SP"
$ws.Range("C4").Value = "ProductId=1"
$ws.Range("D4").Value = "findByProductId(1)"
$ws.Range("E4").Value = "List size > 0 (và có join)"

# --- Row 5 (new): PV_DAO_06 - "Loi Khoa ngoai (Insert)" ---
$ws.Range("A4:G4").Copy()
$ws.Range("A5:G5").PasteSpecial(-4122)
$ws.Range("A5").Value = "PV_DAO_06"
$ws.Range("B5").Value = "Lỗi Khóa ngoại (Insert)"
$ws.Range("C5").Value = "PID=-1"
$ws.Range("D5").Value = "Insert với ProductId rác (-1)"
$ws.Range("E5").Value = "DAO trả về False"
$ws.Range("F5").Value = "OK"
$ws.Range("G5").Value = "PASS"

# --- Row 6 (new): PV_DAO_01 - "Them bien the moi" ---
$ws.Range("A4:G4").Copy()
$ws.Range("A6:G6").PasteSpecial(-4122)
$ws.Range("A6").Value = "PV_DAO_01"
$ws.Range("B6").Value = "Thêm biến thể mới"
$ws.Range("C6").Value = "Stock=99"
$ws.Range("D6").Value = "Insert ProductId=1, Stock=99"
$ws.Range("E6").Value = "Return true"
$ws.Range("F6").Value = "OK"
$ws.Range("G6").Value = "PASS"

# --- Row 7 (new): PV_DAO_05 - "Xoa bien the" ---
$ws.Range("A4:G4").Copy()
$ws.Range("A7:G7").PasteSpecial(-4122)
$ws.Range("A7").Value = "PV_DAO_05"
$ws.Range("B7").Value = "Xóa biến thể"
$ws.Range("C7").Value = "Temp Item"
$ws.Range("D7").Value = "Insert nháp -> Delete"
$ws.Range("E7").Value = "FindById trả về null"
$ws.Range("F7").Value = "OK"
$ws.Range("G7").Value = "PASS"
